# The presentation originally used the "Integral" theme (ppt/theme/theme2.xml,
# referenced by the slide master) while ppt/theme/theme1.xml (referenced only by
# the notes master) held the stock "Office Theme" palette. The authored commit
# swaps the two theme parts' contents, so the deck's active/visible theme
# becomes the stock "Office Theme" color palette (theme1.xml becomes the old
# "Integral" content). The font scheme and format scheme are byte-identical
# between the two themes, so only the 12 color-scheme slots (plus the theme /
# color-scheme display names) actually change.

function ConvertTo-OleColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette: stock PowerPoint "Office Theme" colors, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order (ThemeColorScheme
# item 1..12).
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$theme = $p.Designs.Item(1).SlideMaster.Theme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $theme.ThemeColorScheme.Item($i).RGB = ConvertTo-OleColor $officeThemeHex[$i - 1]
}
